# Atualizacao de bases das ligas, do dia: 04-04-2024 as 23:22
#
# The match rows were re-sorted; for a handful of rows the data (columns B
# through AC) moved to a different row while the row index in column A
# stayed put. This script re-creates that movement by snapshotting the
# "source" row values first and then writing them into the "destination"
# rows, so that every row ends up holding the data that, after the
# re-sort, belongs there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..AC carry the actual match data (column A is just the
# positional row index and does not move).
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

# Maps destination row -> source row (i.e. "the data that ends up in this
# row used to live in the source row").
$rowMap = @{
    122 = 123
    123 = 122
    124 = 125
    125 = 124
    134 = 136
    136 = 135
    135 = 137
    137 = 134
}

# 1) Snapshot every distinct row that is involved, BEFORE any writes, so
#    that overlapping cycles (134 -> 136 -> 135 -> 137 -> 134) don't clobber
#    data we still need to read.
$snapshot = @{}
foreach ($row in $rowMap.Keys) {
    $rowValues = @{}
    foreach ($c in $cols) {
        $addr = $c + $row
        $rowValues[$c] = $ws.Range($addr).Value()
    }
    $snapshot[$row] = $rowValues
}

# 2) Write the snapshotted source-row data into each destination row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcValues = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $addr = $c + $destRow
        $ws.Range($addr).Value = $srcValues[$c]
    }
}
